$d = $word.ActiveDocument
$d.Content.Find.Execute("RX.0X", $true, $false, $false, $false, $false,
                         $true, 1, $false, "RX.XX", 2)
